$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C2 from "Sevilla Ruta" to "Ruta" (reuses existing shared string,
# and the now-unused "Sevilla Ruta" string gets dropped on save)
$ws.Range("C2").Value = "Ruta"

# Reflect the new active selection on the sheet (as in the authored change)
[void]$ws.Range("C2").Select()
